# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" worksheet (fund-holdings detail, same layout as
#   the other quarterly sheets) right before the "总计" (totals) sheet.
# - Prepend a new row to the "总计" sheet summarising the 2022-Q1 quarter,
#   pushing the existing history rows down by one and renumbering the
#   index column.

$wb = $excel.ActiveWorkbook

# A sheet used as a formatting donor for the bold/centered/bordered header
# style (cellXfs index 2 in the original file) that every quarterly sheet
# shares for its header row and its "A" index column.
$donor = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet immediately before "总计".
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

# Copy header / index-column formatting from the donor sheet so the new
# sheet matches the look of the existing quarterly sheets.
$donor.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$donor.Range("A2").Copy()
$newSheet.Range("A2:A15").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Header row.
# ---------------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# 3. Fund holdings data rows (A..H), rows 2-15.
# ---------------------------------------------------------------------------
$data = @(
    @(0,  "161903", "万家行业优选混合 (LOF)", "112.51", "91.94", "7.17", "8.0670", 6),
    @(1,  "161914", "万家创业板2年定期开放混合A", "15.74", "95.20", "8.54", "1.3442", 5),
    @(2,  "005311", "万家经济新动能混合A", "19.47", "93.80", "6.35", "1.2363", 9),
    @(3,  "005312", "万家经济新动能混合C", "6.77", "93.80", "6.35", "0.4299", 9),
    @(4,  "159883", "永赢中证全指医疗器械交易型开放式指数证券投资基金", "9.17", "99.23", "3.00", "0.2751", 7),
    @(5,  "161915", "万家创业板2年定期开放混合C", "2.36", "95.20", "8.54", "0.2015", 5),
    @(6,  "159898", "招商中证全指医疗器械交易型开放式指数证券投资基金", "1.61", "99.41", "2.96", "0.0477", 7),
    @(7,  "167506", "安信中证深圳科技创新主题指数（LOF）A", "1.27", "90.25", "3.61", "0.0458", 8),
    @(8,  "159873", "天弘中证全指医疗保健设备与服务ETF", "1.43", "99.59", "2.57", "0.0368", 9),
    @(9,  "159891", "建信中证全指医疗保健设备与服务交易型开放式指数证券投资基金", "1.40", "95.24", "2.51", "0.0351", 9),
    @(10, "516610", "大成中证全指医疗保健设备与服务交易型开放式指数证券投资基金", "0.69", "96.20", "2.36", "0.0163", 9),
    @(11, "167507", "安信中证深圳科技创新主题指数（LOF）C", "0.33", "90.25", "3.61", "0.0119", 8),
    @(12, "501069", "华宝标普中国Ａ股质量价值指数（ＬＯＦ）", "0.16", "94.73", "2.97", "0.0048", 4),
    @(13, "001797", "华融新利灵活配置混合", "0.02", "48.66", "2.29", "0.0005", 8)
)

$r = 2
foreach ($row in $data) {
    # Columns B, D, E, F, G hold numeric-looking values that must stay text
    # (so fund codes keep leading zeros and decimals keep trailing zeros,
    # e.g. "005311" / "8.0670") - a leading apostrophe forces text entry
    # exactly the way typing it into Excel would.
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 4. Update the "总计" sheet: insert a new top data row for 2022-Q1 and
#    shift the previously-existing history down by one row.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Re-apply the data-row formatting (A column index style) that Insert()
# doesn't carry over correctly, by copying it from the row just below.
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 14
$totalSheet.Cells.Item(2, 4).Value = 11.75

# Renumber the 0-based index column for the rows that got pushed down.
for ($i = 3; $i -le 7; $i++) {
    $totalSheet.Cells.Item($i, 1).Value = $i - 2
}
